$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.618
$ws.Range("A3").Value = -21.557
$ws.Range("D5").Value = -8.128
$ws.Range("A14").Value = -21.04
$ws.Range("A16").Value = -20.771
$ws.Range("D16").Value = -8.488
$ws.Range("A21").Value = -21.04
$ws.Range("A23").Value = -21.584
$ws.Range("A25").Value = -22.078
